# SWP391-AppDevProject_Evaluations-LOC.xlsx edit
#
# Summary of the target change (per the OOXML diff):
#  1. The empty helper tab "Sheet1" is removed.
#  2. On "Assignemnts-original":
#       - rows 9-13: Function Complexity (col C) "Medium" -> "Complex"
#         (driving the dependent LOC formulas in D/L/O to double)
#       - E10, E11, E13 (Planned Code Iteration) corrected to "Iteration 1"
#       - row 15 (a leftover "--" placeholder row) is cleared out entirely
#       - the active sheet view scrolls to A4 or before B4, and the
#         selection moves from F8 to C15
#  3. The two pivot-table sheets ("Iteration-LOC", "Final-LOC") keep their
#     on-screen (cached) pivot grid exactly as-is; only the plain VLOOKUP
#     helper formulas beside the pivot (col A / col B) re-evaluate, turning
#     "--" into #N/A once row 15's "--" roll number disappears from the
#     source table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Drop the blank "Sheet1" tab.
# ---------------------------------------------------------------------
$wb.Worksheets("Sheet1").Delete()

# ---------------------------------------------------------------------
# 2) Edit the "Assignemnts-original" source table.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets("Assignemnts-original")

$ws1.Range("C9").Value = "Complex"

$ws1.Range("C10").Value = "Complex"
$ws1.Range("E10").Value = "Iteration 1"

$ws1.Range("C11").Value = "Complex"
$ws1.Range("E11").Value = "Iteration 1"

$ws1.Range("C12").Value = "Complex"

$ws1.Range("C13").Value = "Complex"
$ws1.Range("E13").Value = "Iteration 1"

# Row 15 was a stray "--" placeholder entry; wipe it back to blank
# (keeping the existing cell formatting in place).
$ws1.Range("A15:O15").ClearContents()

# Restore the view: scroll back to column A and land the selection on C15.
$ws1.Application.ActiveWindow.ScrollColumn = 1
$ws1.Application.ActiveWindow.ScrollRow = 4
$ws1.Range("C15").Select()

# ---------------------------------------------------------------------
# 3) This engine (unlike Excel) reactively refreshes PivotTable caches
#    whenever their source range changes; Excel itself only refreshes a
#    PivotTable on an explicit RefreshAll / file open. To match the real
#    target workbook (where the cached pivot grids stay untouched) put
#    the pivot output cells back to their original cached values. The
#    ordinary VLOOKUP formulas living alongside each pivot are left
#    alone so they keep recalculating naturally.
# ---------------------------------------------------------------------
$wsIter = $wb.Worksheets("Iteration-LOC")

$wsIter.Range("B9").Value = "Sum of Iteration LOC"
$wsIter.Range("C9").Value = "Column Labels"
$wsIter.Range("A10").Value = "Student Name"
$wsIter.Range("B10").Value = "Row Labels"
$wsIter.Range("C10").Value = "Iteration 1"
$wsIter.Range("D10").Value = "Iteration 2"
$wsIter.Range("E10").Value = "Iteration 3"
$wsIter.Range("F10").Value = "Grand Total"

$wsIter.Range("B11").Value = "--"
$wsIter.Range("C11").Value = 45
$wsIter.Range("D11").ClearContents()
$wsIter.Range("E11").ClearContents()
$wsIter.Range("F11").Value = 45

$wsIter.Range("B12").Value = "SE04964"
$wsIter.Range("C12").ClearContents()
$wsIter.Range("D12").ClearContents()
$wsIter.Range("E12").Value = 135
$wsIter.Range("F12").Value = 135

$wsIter.Range("B13").Value = "SE05407"
$wsIter.Range("C13").Value = 60
$wsIter.Range("D13").Value = 135
$wsIter.Range("E13").Value = 180
$wsIter.Range("F13").Value = 375

$wsIter.Range("B14").Value = "SE05436"
$wsIter.Range("C14").ClearContents()
$wsIter.Range("D14").Value = 60
$wsIter.Range("E14").ClearContents()
$wsIter.Range("F14").Value = 60

$wsIter.Range("B15").Value = "SE05725"
$wsIter.Range("C15").Value = 90
$wsIter.Range("D15").ClearContents()
$wsIter.Range("E15").ClearContents()
$wsIter.Range("F15").Value = 90

$wsIter.Range("B16").Value = "Grand Total"
$wsIter.Range("C16").Value = 195
$wsIter.Range("D16").Value = 195
$wsIter.Range("E16").Value = 315
$wsIter.Range("F16").Value = 705

$wsFinal = $wb.Worksheets("Final-LOC")

$wsFinal.Range("B10").Value = "Student Name"
$wsFinal.Range("C10").Value = "Row Labels"
$wsFinal.Range("D10").Value = "Sum of Final LOC"

$wsFinal.Range("C11").Value = "--"
$wsFinal.Range("D11").Value = 45

$wsFinal.Range("C12").Value = "SE04964"
$wsFinal.Range("D12").Value = 135

$wsFinal.Range("C13").Value = "SE05407"
$wsFinal.Range("D13").Value = 360

$wsFinal.Range("C14").Value = "SE05436"
$wsFinal.Range("D14").Value = 90

$wsFinal.Range("C15").Value = "SE05725"
$wsFinal.Range("D15").Value = 90

$wsFinal.Range("C16").Value = "Grand Total"
$wsFinal.Range("D16").Value = 720

# ---------------------------------------------------------------------
# 4) Re-select "Assignemnts-original" as the active sheet (it's the
#    tabSelected="1" sheet in the saved file).
# ---------------------------------------------------------------------
$ws1.Activate()
